$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two raw input cells that drive the rest of the model via formulas.
$ws.Range("D15").Value = 8
$ws.Range("D17").Value = 7.5

# Reflect the cell selected on Sheet1 after making the edit.
$ws.Range("D12").Select()
